$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "jorge_jorge"
$ws.Range("C3").Value = "123456Senha"
$ws.Range("D3").Value = "Cadastrado"

$ws.Range("B5").Value = "jorge_egorj"
$ws.Range("C5").Value = "123456Senha"
$ws.Range("D5").Value = "UserName Incorreto"

$ws.Range("B3:D3").HorizontalAlignment = -4108
$ws.Range("B3:D3").VerticalAlignment = -4108
$ws.Range("B5:D5").HorizontalAlignment = -4108
$ws.Range("B5:D5").VerticalAlignment = -4108

$ws.Columns.Item(4).ColumnWidth = 19.85546875

$ws.Cells.Item(11, 5).Select()
